# Applies the cryptos list price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => { column letter => new cell text }
$updates = @{
    2 = @{ D="44.419.09"; E="  +0.55%  " }
    3 = @{ D="2.245.18"; E="  -0.42%  " }
    4 = @{ E="  +0.44%  " }
    5 = @{ D="305.97"; E="  -0.49%  " }
    6 = @{ D="93.15"; E="  -6.06%  " }
    7 = @{ D="0.570"; E="  -1.00%  " }
    9 = @{ D="0.522"; E="  -2.73%  " }
    10 = @{ D="34.53"; E="  -3.12%  " }
    11 = @{ E="  -1.66%  " }
    12 = @{ D="7.12"; E="  -2.99%  " }
    13 = @{ E="  -0.21%  " }
    14 = @{ D="2.250.37"; E="  -0.03%  " }
    15 = @{ D="0.836"; E="  -0.69%  " }
    16 = @{ D="13.54"; E="  -2.37%  " }
    17 = @{ D="44.110.47"; E="  +0.13%  " }
    18 = @{ D="0.0₃0960"; E="  -1.90%  " }
    19 = @{ D="12.31"; E="  -4.46%  " }
    20 = @{ E="  -0.17%  " }
    21 = @{ D="65.51"; E="  -0.01%  " }
    22 = @{ B="PancakeSwap"; C="https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D="3.13"; E="  +5.97%  " }
    23 = @{ B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="237.31"; E="  -2.02%  " }
    24 = @{ E="  -0.10%  " }
    25 = @{ E="  -0.21%  " }
    26 = @{ D="38.72"; E="  +4.52%  " }
    27 = @{ E="  +2.92%  " }
    28 = @{ D="9.79"; E="  -3.44%  " }
    29 = @{ D="5.91"; E="  -3.84%  " }
    30 = @{ D="20.04"; E="  -0.60%  " }
    31 = @{ D="153.85"; E="  -2.56%  " }
    32 = @{ D="0.0796"; E="  -4.11%  " }
    33 = @{ D="2.65"; E="  -0.30%  " }
    34 = @{ E="  -14.32%  " }
    35 = @{ E="  +0.25%  " }
    36 = @{ E="  -0.38%  " }
    37 = @{ D="1.82"; E="  -2.37%  " }
    38 = @{ D="3.44"; E="  +1.12%  " }
    39 = @{ D="14.53"; E="  -4.87%  " }
    40 = @{ D="3.81"; E="  -2.21%  " }
    41 = @{ D="0.0301"; E="  -2.12%  " }
    42 = @{ E="  +0.36%  " }
    43 = @{ D="1.731.62"; E="  -1.48%  " }
    44 = @{ B="Algorand"; C="https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; D="0.192"; E="  -0.53%  " }
    45 = @{ B="BitcoinSV"; C="https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"; D="80.24"; E="  -9.71%  " }
    46 = @{ D="99.28"; E="  -2.35%  " }
    47 = @{ E="  -4.66%  " }
    48 = @{ E="  +3.66%  " }
    49 = @{ D="8.17"; E="  -1.19%  " }
    50 = @{ B="MultiversX"; C="https://coinranking.com/coin/omwkOTglq+multiversx-egld"; D="55.22"; E="  -0.63%  " }
    51 = @{ B="ordi"; C="https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"; D="69.40"; E="  -1.46%  " }
}

# Cells whose new text would otherwise be auto-parsed as a number by Excel,
# losing a significant trailing zero (e.g. "0.570" -> 0.57). Force Text format
# on those specific cells before writing the value so they round-trip exactly.
$forceTextCells = @("D7", "D51")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

foreach ($r in $updates.Keys) {
    $cols = $updates[$r]
    foreach ($col in $cols.Keys) {
        $addr = "$col$r"
        $ws.Range($addr).Value = $cols[$col]
    }
}

Write-Host "Updated cryptos list ($($updates.Keys.Count) rows)."